$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 347) holds the "Förändrad" date value.
# All of these cells currently store the serial date 45172 and should
# be updated to the new serial date 45175.
$ws.Range("C2:C347").Value = 45175
